{"js": "// Apply the descr-table refresh: categorical row counts/percentages were\n// recomputed (now an even 15/15 split instead of 14/18, 16/12), the\n// \"a\" row's chi-squared p-value and Wald CI were updated, and the\n// narrower \"p\" column widened slightly to fit the new value.\n\n// 1) Widen the \"p\" column (4th gridCol, 0-based index 4) from 1212 -> 1291\n//    twips. Word expresses cell/column widths in points, so divide by 20.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst targetTable = tables.items[3];\nconst pColumnHeaderCell = targetTable.getCell(0, 4);\npColumnHeaderCell.columnWidth = 1291 / 20;\nawait context.sync();\n\n// 2) Replace the stale summary-statistic text runs. Each value is a\n//    unique string in the document, so a plain text search+replace keeps\n//    the surrounding run/paragraph formatting untouched.\nconst replacements = [\n  [\"14 (47%)\", \"15 (50%)\"],\n  [\"18 (60%)\", \"15 (50%)\"],\n  [\"32 (53%)\", \"30 (50%)\"],\n  [\"0.301\", \">0.999\"],\n  [\"[-0.39, 0.12]\", \"[-0.25, 0.25]\"],\n  [\"16 (53%)\", \"15 (50%)\"],\n  [\"12 (40%)\", \"15 (50%)\"],\n  [\"28 (47%)\", \"30 (50%)\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the descr-table refresh: categorical row counts/percentages were\n# recomputed (now an even 15/15 split instead of 14/18, 16/12), the\n# \"a\" row's chi-squared p-value and Wald CI were updated, and the\n# narrower \"p\" column widened slightly to fit the new value.\n\n$d = $word.ActiveDocument\n\n# 1) Widen the \"p\" column (5th column, 1-based) from 1212 -> 1291 twips.\n#    Word expresses column width in points, so divide by 20.\n$t = $d.Tables.Item(4)\n$t.Columns.Item(5).Width = 1291 / 20\n\n# 2) Replace the stale summary-statistic text runs. Each value is a\n#    unique string in the document, so a plain Find/Replace keeps the\n#    surrounding run/paragraph formatting untouched.\n$replacements = @(\n    @(\"14 (47%)\", \"15 (50%)\"),\n    @(\"18 (60%)\", \"15 (50%)\"),\n    @(\"32 (53%)\", \"30 (50%)\"),\n    @(\"0.301\", \">0.999\"),\n    @(\"[-0.39, 0.12]\", \"[-0.25, 0.25]\"),\n    @(\"16 (53%)\", \"15 (50%)\"),\n    @(\"12 (40%)\", \"15 (50%)\"),\n    @(\"28 (47%)\", \"30 (50%)\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
